$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LTE001_ACC_00001")

# New data rows to append (row 2 already exists in the sheet)
$newRows = @(
    @{ A=11377; B=11377; C=11377; D="SEA"; E="LAX"; F="GOLDSTREAK"; G="None"; H="NONSCR"; I="None"; J="None"; K=20; L=750; M="PP"; N="CREDIT"; O="CART" },
    @{ A=11377; B=11377; C=11377; D="DFW"; E="SEA"; F="GENERAL";    G="None"; H="NONSCR"; I="None"; J="None"; K=1;  L=100; M="PP"; N="CREDIT"; O="CART" },
    @{ A=11377; B=11377; C=11377; D="DFW"; E="SEA"; F="PRIORITY";   G="None"; H="NONSCR"; I="None"; J="None"; K=8;  L=600; M="PP"; N="CREDIT"; O="CART" },
    @{ A=11377; B=11377; C=11377; D="SEA"; E="ANC"; F="GENERAL";    G="None"; H="NONSCR"; I="None"; J="None"; K=1;  L=100; M="PP"; N="CREDIT"; O="CART" }
)

$rowIndex = 3
foreach ($rowData in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData.A
    $ws.Cells.Item($rowIndex, 2).Value = $rowData.B
    $ws.Cells.Item($rowIndex, 3).Value = $rowData.C
    $ws.Cells.Item($rowIndex, 4).Value = $rowData.D
    $ws.Cells.Item($rowIndex, 5).Value = $rowData.E
    $ws.Cells.Item($rowIndex, 6).Value = $rowData.F
    $ws.Cells.Item($rowIndex, 7).Value = $rowData.G
    $ws.Cells.Item($rowIndex, 8).Value = $rowData.H
    $ws.Cells.Item($rowIndex, 9).Value = $rowData.I
    $ws.Cells.Item($rowIndex, 10).Value = $rowData.J
    $ws.Cells.Item($rowIndex, 11).Value = $rowData.K
    $ws.Cells.Item($rowIndex, 12).Value = $rowData.L
    $ws.Cells.Item($rowIndex, 13).Value = $rowData.M
    $ws.Cells.Item($rowIndex, 14).Value = $rowData.N
    $ws.Cells.Item($rowIndex, 15).Value = $rowData.O
    $rowIndex++
}

# Make LTE001_ACC_00001 the active sheet/tab and set its selection
$ws.Activate()
$ws.Range("F6").Select() | Out-Null
